$wb = $excel.ActiveWorkbook

# Sheet ALC, row 7
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 1158.7273
$ws.Range("J7").Value = 2500
$ws.Range("L7").Value = 2500
$ws.Range("N7").Value = -2724

# Sheet ALC, row 14
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H14").Value = 1158.7273
$ws.Range("J14").Value = 2500
$ws.Range("L14").Value = 2500
$ws.Range("N14").Value = -2882

# Sheet ALC, row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3658
$ws.Range("J18").Value = 6316
$ws.Range("L18").Value = 6316
$ws.Range("N18").Value = -6884

# Sheet ALC, row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1333.3334
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 2000
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 2000
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -2138

# Sheet ALC, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7211.875
$ws.Range("I62").Value = 3423.75
$ws.Range("J62").Value = 11000
$ws.Range("K62").Value = 3423.75
$ws.Range("L62").Value = 11000
$ws.Range("M62").Value = -2799.75
$ws.Range("N62").Value = -12248

# Sheet ALC, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 7211.875
$ws.Range("I65").Value = 3423.75
$ws.Range("J65").Value = 11000
$ws.Range("K65").Value = 17118.75
$ws.Range("L65").Value = 55000
$ws.Range("M65").Value = -13998.75
$ws.Range("N65").Value = -61240

# Sheet ALC, row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1001
$ws.Range("I107").Value = 999
$ws.Range("J107").Value = 1002
$ws.Range("K107").Value = 999
$ws.Range("L107").Value = 1002
$ws.Range("M107").Value = 921
$ws.Range("N107").Value = -4842

# Sheet ALC, row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1972.5
$ws.Range("I111").Value = 1300
$ws.Range("J111").Value = 3990
$ws.Range("K111").Value = 3900
$ws.Range("L111").Value = 11970
$ws.Range("M111").Value = -833
$ws.Range("N111").Value = -18104

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2051.6
$ws.Range("I137").Value = 1576.2
$ws.Range("J137").Value = 3002.4
$ws.Range("K137").Value = 4728.6
$ws.Range("L137").Value = 9007.200000000001
$ws.Range("M137").Value = -2178.6
$ws.Range("N137").Value = -14107.2

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2362.1667
$ws.Range("I141").Value = 2362.1667
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 7086.500100000001
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -1906.500100000001
$ws.Range("N141").ClearContents()

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2416.6667
$ws.Range("I2").Value = 2416.6667
$ws.Range("K2").Value = 2416.6667
$ws.Range("M2").Value = -2303.6667

# Sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4933.3335
$ws.Range("I45").Value = 4933.3335
$ws.Range("K45").Value = 4933.3335
$ws.Range("M45").Value = -4556.3335

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1968
$ws.Range("I61").Value = 1968
$ws.Range("K61").Value = 1968
$ws.Range("M61").Value = -1756

# Sheet ARM, row 76
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 85000
$ws.Range("J76").Value = 85000
$ws.Range("L76").Value = 85000
$ws.Range("N76").Value = -85676

# Sheet ARM, row 79
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 85000
$ws.Range("J79").Value = 85000
$ws.Range("L79").Value = 85000
$ws.Range("N79").Value = -87340

# Sheet ARM, row 92
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 44880
$ws.Range("J92").Value = 44880
$ws.Range("L92").Value = 44880
$ws.Range("N92").Value = -49872

# Sheet ARM, row 95
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 13517.143
$ws.Range("J95").Value = 13517.143
$ws.Range("L95").Value = 13517.143
$ws.Range("N95").Value = -19009.143

# Sheet ARM, row 96
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 33355.25
$ws.Range("J96").Value = 33355.25
$ws.Range("L96").Value = 33355.25
$ws.Range("N96").Value = -38847.25

# Sheet ARM, row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 656.1111
$ws.Range("I97").Value = 515
$ws.Range("J97").Value = 1150
$ws.Range("K97").Value = 515
$ws.Range("L97").Value = 1150
$ws.Range("M97").Value = -19
$ws.Range("N97").Value = -2142

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2416.6667
$ws.Range("I116").Value = 2416.6667
$ws.Range("K116").Value = 2416.6667
$ws.Range("M116").Value = -122.6667000000002

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I122").Value = 1381
$ws.Range("J122").Value = 1490
$ws.Range("K122").Value = 4143
$ws.Range("L122").Value = 4470
$ws.Range("M122").Value = -1693
$ws.Range("N122").Value = -9370

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2910.8823
$ws.Range("I132").Value = 1649.8462
$ws.Range("K132").Value = 4949.5386
$ws.Range("M132").Value = -2419.5386

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1968
$ws.Range("I136").Value = 1968
$ws.Range("K136").Value = 5904
$ws.Range("M136").Value = -3354

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2416.6667
$ws.Range("I3").Value = 2416.6667
$ws.Range("K3").Value = 2416.6667
$ws.Range("M3").Value = -2302.6667

# Sheet BSM, row 92
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 40000.5
$ws.Range("J92").Value = 40000.5
$ws.Range("L92").Value = 40000.5
$ws.Range("N92").Value = -44992.5

# Sheet CUL, row 6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 15.5
$ws.Range("I6").Value = 16.333334
$ws.Range("K6").Value = 49.000002
$ws.Range("M6").Value = 63.999998

# Sheet CUL, row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 567.45
$ws.Range("J68").Value = 543.6
$ws.Range("L68").Value = 1630.8
$ws.Range("N68").Value = -3252.8

# Sheet CUL, row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 567.45
$ws.Range("J71").Value = 543.6
$ws.Range("L71").Value = 4892.400000000001
$ws.Range("N71").Value = -13004.4

# Sheet CUL, row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 316.5
$ws.Range("I122").Value = 316.5
$ws.Range("K122").Value = 2848.5
$ws.Range("M122").Value = -398.5

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("N131").ClearContents()

# Sheet CUL, row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 4266.5
$ws.Range("I140").Value = 4500
$ws.Range("K140").Value = 13500
$ws.Range("M140").Value = -8320

# Sheet GSM, row 3
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1695495.1
$ws.Range("I3").Value = 1314001.5
$ws.Range("J3").Value = 2000690
$ws.Range("K3").Value = 1314001.5
$ws.Range("L3").Value = 2000690
$ws.Range("M3").Value = -1313885.5
$ws.Range("N3").Value = -2000922

# Sheet GSM, row 9
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()

# Sheet GSM, row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8350.546
$ws.Range("J70").Value = 9099.75
$ws.Range("L70").Value = 9099.75
$ws.Range("N70").Value = -9639.75

# Sheet GSM, row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 8350.546
$ws.Range("J73").Value = 9099.75
$ws.Range("L73").Value = 9099.75
$ws.Range("N73").Value = -10971.75

# Sheet GSM, row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 857.3333
$ws.Range("I107").Value = 832
$ws.Range("K107").Value = 832
$ws.Range("M107").Value = 1088

# Sheet LTW, row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1847.5
$ws.Range("I55").Value = 695
$ws.Range("J55").Value = 3000
$ws.Range("K55").Value = 695
$ws.Range("L55").Value = 3000
$ws.Range("M55").Value = -522
$ws.Range("N55").Value = -3346

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9748.75
$ws.Range("I122").Value = 9747.5
$ws.Range("K122").Value = 29242.5
$ws.Range("M122").Value = -26792.5

# Sheet WVR, row 17
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 252.5
$ws.Range("I17").Value = 252.5
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 252.5
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -80.5
$ws.Range("N17").ClearContents()

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2360.5
$ws.Range("I122").Value = 1471
$ws.Range("K122").Value = 4413
$ws.Range("M122").Value = -1963

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3000
$ws.Range("I136").Value = 1250
$ws.Range("K136").Value = 3750
$ws.Range("M136").Value = -1200
